# Update odds figures (Betfair Back/Lay) for 2026-01-20 fixtures.
# Only numeric data cells (columns F:AO) change; League/Date/Time/Home/Away
# (columns A:E) and the header row are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 6).Value = 1.55
$ws.Cells.Item(2, 7).Value = 1.64
$ws.Cells.Item(2, 8).Value = 6
$ws.Cells.Item(2, 9).Value = 7.8
$ws.Cells.Item(2, 10).Value = 4.1
$ws.Cells.Item(2, 14).Value = 3.85
$ws.Cells.Item(2, 20).Value = 1.91
$ws.Cells.Item(2, 23).Value = 2.56
$ws.Cells.Item(2, 34).Value = 24
$ws.Cells.Item(2, 37).Value = 17.5

# Row 3
$ws.Cells.Item(3, 7).Value = 5.8
$ws.Cells.Item(3, 8).Value = 1.61
$ws.Cells.Item(3, 9).Value = 1.63
$ws.Cells.Item(3, 12).Value = 1.32
$ws.Cells.Item(3, 14).Value = 5.2
$ws.Cells.Item(3, 15).Value = 1.22
$ws.Cells.Item(3, 16).Value = 2.38
$ws.Cells.Item(3, 17).Value = 1.68
$ws.Cells.Item(3, 18).Value = 1.55
$ws.Cells.Item(3, 19).Value = 2.74
$ws.Cells.Item(3, 20).Value = 1.75
$ws.Cells.Item(3, 21).Value = 2.22
$ws.Cells.Item(3, 22).Value = 2.58
$ws.Cells.Item(3, 26).Value = 11
$ws.Cells.Item(3, 28).Value = 25
$ws.Cells.Item(3, 30).Value = 9.8
$ws.Cells.Item(3, 32).Value = 48
$ws.Cells.Item(3, 35).Value = 29
$ws.Cells.Item(3, 36).Value = 150
$ws.Cells.Item(3, 37).Value = 70
$ws.Cells.Item(3, 38).Value = 65
$ws.Cells.Item(3, 39).Value = 90
$ws.Cells.Item(3, 41).Value = 7.4

# Row 4
$ws.Cells.Item(4, 6).Value = 2.12
$ws.Cells.Item(4, 8).Value = 3.25
$ws.Cells.Item(4, 9).Value = 3.8
$ws.Cells.Item(4, 10).Value = 3.15
$ws.Cells.Item(4, 12).Value = 1.01
$ws.Cells.Item(4, 13).Value = 1.06
$ws.Cells.Item(4, 14).Value = 1.98
$ws.Cells.Item(4, 15).Value = 1.27
$ws.Cells.Item(4, 18).Value = 1.18
$ws.Cells.Item(4, 19).Value = 1.83
$ws.Cells.Item(4, 20).Value = 1.01
$ws.Cells.Item(4, 21).Value = 1.01
$ws.Cells.Item(4, 22).Value = 1.35
$ws.Cells.Item(4, 23).Value = 1.73
$ws.Cells.Item(4, 24).Value = 1000
$ws.Cells.Item(4, 25).Value = 20
$ws.Cells.Item(4, 26).Value = 36
$ws.Cells.Item(4, 27).Value = 1000
$ws.Cells.Item(4, 28).Value = 15
$ws.Cells.Item(4, 29).Value = 11.5
$ws.Cells.Item(4, 30).Value = 20
$ws.Cells.Item(4, 31).Value = 1000
$ws.Cells.Item(4, 32).Value = 21
$ws.Cells.Item(4, 33).Value = 15.5
$ws.Cells.Item(4, 34).Value = 23
$ws.Cells.Item(4, 35).Value = 1000
$ws.Cells.Item(4, 36).Value = 40
$ws.Cells.Item(4, 37).Value = 32
$ws.Cells.Item(4, 38).Value = 1000
$ws.Cells.Item(4, 39).Value = 1000
$ws.Cells.Item(4, 40).Value = 1000
$ws.Cells.Item(4, 41).Value = 1000

# Row 5
$ws.Cells.Item(5, 6).Value = 1.7
$ws.Cells.Item(5, 7).Value = 1.78
$ws.Cells.Item(5, 9).Value = 5.3
$ws.Cells.Item(5, 11).Value = 4.9
$ws.Cells.Item(5, 12).Value = 1.01
$ws.Cells.Item(5, 13).Value = 1.03
$ws.Cells.Item(5, 14).Value = 5.4
$ws.Cells.Item(5, 15).Value = 1.14
$ws.Cells.Item(5, 16).Value = 2.68
$ws.Cells.Item(5, 17).Value = 1.48
$ws.Cells.Item(5, 18).Value = 1.68
$ws.Cells.Item(5, 19).Value = 2.18
$ws.Cells.Item(5, 20).Value = 1.55
$ws.Cells.Item(5, 21).Value = 2.46
$ws.Cells.Item(5, 22).Value = 1.25
$ws.Cells.Item(5, 23).Value = 2.28
$ws.Cells.Item(5, 24).Value = 34
$ws.Cells.Item(5, 25).Value = 30
$ws.Cells.Item(5, 26).Value = 50
$ws.Cells.Item(5, 27).Value = 120
$ws.Cells.Item(5, 28).Value = 16
$ws.Cells.Item(5, 29).Value = 13.5
$ws.Cells.Item(5, 30).Value = 23
$ws.Cells.Item(5, 31).Value = 60
$ws.Cells.Item(5, 32).Value = 16
$ws.Cells.Item(5, 33).Value = 12.5
$ws.Cells.Item(5, 34).Value = 19.5
$ws.Cells.Item(5, 35).Value = 55
$ws.Cells.Item(5, 36).Value = 22
$ws.Cells.Item(5, 37).Value = 18.5
$ws.Cells.Item(5, 38).Value = 29
$ws.Cells.Item(5, 39).Value = 70
$ws.Cells.Item(5, 40).Value = 7.4
$ws.Cells.Item(5, 41).Value = 42

# Row 6
$ws.Cells.Item(6, 7).Value = 7.6
$ws.Cells.Item(6, 9).Value = 1.43
$ws.Cells.Item(6, 10).Value = 6.2
$ws.Cells.Item(6, 11).Value = 6.4
$ws.Cells.Item(6, 12).Value = 1.19
$ws.Cells.Item(6, 14).Value = 8
$ws.Cells.Item(6, 16).Value = 3.4
$ws.Cells.Item(6, 18).Value = 1.98
$ws.Cells.Item(6, 19).Value = 1.95
$ws.Cells.Item(6, 20).Value = 1.62
$ws.Cells.Item(6, 22).Value = 3.3
$ws.Cells.Item(6, 23).Value = 1.15
$ws.Cells.Item(6, 24).Value = 42
$ws.Cells.Item(6, 25).Value = 15.5
$ws.Cells.Item(6, 26).Value = 12.5
$ws.Cells.Item(6, 27).Value = 14
$ws.Cells.Item(6, 28).Value = 42
$ws.Cells.Item(6, 29).Value = 14.5
$ws.Cells.Item(6, 30).Value = 10.5
$ws.Cells.Item(6, 31).Value = 12.5
$ws.Cells.Item(6, 32).Value = 75
$ws.Cells.Item(6, 33).Value = 29
$ws.Cells.Item(6, 34).Value = 20
$ws.Cells.Item(6, 35).Value = 24
$ws.Cells.Item(6, 36).Value = 220
$ws.Cells.Item(6, 37).Value = 80
$ws.Cells.Item(6, 38).Value = 65
$ws.Cells.Item(6, 39).Value = 70
$ws.Cells.Item(6, 40).Value = 55
$ws.Cells.Item(6, 41).Value = 3.95

# Row 7
$ws.Cells.Item(7, 17).Value = 2.12

# Row 8
$ws.Cells.Item(8, 6).Value = 1.6
$ws.Cells.Item(8, 9).Value = 6.4
$ws.Cells.Item(8, 11).Value = 4.8

# Row 10
$ws.Cells.Item(10, 6).Value = 3
$ws.Cells.Item(10, 7).Value = 3.45
$ws.Cells.Item(10, 8).Value = 2.28
$ws.Cells.Item(10, 9).Value = 2.46
$ws.Cells.Item(10, 10).Value = 3.65
$ws.Cells.Item(10, 11).Value = 3.95
$ws.Cells.Item(10, 16).Value = 2.06
$ws.Cells.Item(10, 17).Value = 1.75

# Row 12
$ws.Cells.Item(12, 6).Value = 1.73
$ws.Cells.Item(12, 8).Value = 3.8
$ws.Cells.Item(12, 11).Value = 6.4

# Row 13
$ws.Cells.Item(13, 17).Value = 1.78
$ws.Cells.Item(13, 19).Value = 2.96
$ws.Cells.Item(13, 27).Value = 48
$ws.Cells.Item(13, 28).Value = 13.5
$ws.Cells.Item(13, 30).Value = 13
$ws.Cells.Item(13, 33).Value = 12.5
$ws.Cells.Item(13, 37).Value = 1000
$ws.Cells.Item(13, 41).Value = 30

# Row 14
$ws.Cells.Item(14, 9).Value = 1.92
$ws.Cells.Item(14, 15).Value = 1.35
$ws.Cells.Item(14, 16).Value = 1.86
$ws.Cells.Item(14, 17).Value = 2.14
$ws.Cells.Item(14, 19).Value = 3.7
$ws.Cells.Item(14, 20).Value = 1.89

# Row 15
$ws.Cells.Item(15, 6).Value = 2.86
$ws.Cells.Item(15, 7).Value = 2.88
$ws.Cells.Item(15, 8).Value = 2.68
$ws.Cells.Item(15, 9).Value = 2.7
$ws.Cells.Item(15, 17).Value = 2.02
$ws.Cells.Item(15, 18).Value = 1.34
$ws.Cells.Item(15, 19).Value = 3.65
$ws.Cells.Item(15, 20).Value = 1.79
$ws.Cells.Item(15, 26).Value = 18
$ws.Cells.Item(15, 32).Value = 20
$ws.Cells.Item(15, 41).Value = 1000

# Row 16
$ws.Cells.Item(16, 6).Value = 4.4
$ws.Cells.Item(16, 7).Value = 4.7
$ws.Cells.Item(16, 8).Value = 1.78
$ws.Cells.Item(16, 9).Value = 1.81
$ws.Cells.Item(16, 10).Value = 4.3
$ws.Cells.Item(16, 17).Value = 1.64
$ws.Cells.Item(16, 18).Value = 1.57
$ws.Cells.Item(16, 19).Value = 2.62
$ws.Cells.Item(16, 20).Value = 1.65
$ws.Cells.Item(16, 21).Value = 2.46
$ws.Cells.Item(16, 26).Value = 13.5
$ws.Cells.Item(16, 31).Value = 17
$ws.Cells.Item(16, 32).Value = 38
$ws.Cells.Item(16, 34).Value = 17.5
$ws.Cells.Item(16, 36).Value = 1000
$ws.Cells.Item(16, 37).Value = 60
$ws.Cells.Item(16, 38).Value = 60

# Row 17
$ws.Cells.Item(17, 11).Value = 5
$ws.Cells.Item(17, 16).Value = 2.92
$ws.Cells.Item(17, 29).Value = 12
$ws.Cells.Item(17, 37).Value = 15

# Row 18
$ws.Cells.Item(18, 16).Value = 3.35
$ws.Cells.Item(18, 18).Value = 1.97
$ws.Cells.Item(18, 20).Value = 1.74
$ws.Cells.Item(18, 21).Value = 2.26
$ws.Cells.Item(18, 32).Value = 11.5
$ws.Cells.Item(18, 33).Value = 11.5
$ws.Cells.Item(18, 34).Value = 25
$ws.Cells.Item(18, 36).Value = 12
$ws.Cells.Item(18, 37).Value = 12.5
$ws.Cells.Item(18, 38).Value = 29
$ws.Cells.Item(18, 39).Value = 95

# Row 19
$ws.Cells.Item(19, 6).Value = 2.78
$ws.Cells.Item(19, 7).Value = 2.84
$ws.Cells.Item(19, 8).Value = 2.64
$ws.Cells.Item(19, 9).Value = 2.66
$ws.Cells.Item(19, 17).Value = 1.82
$ws.Cells.Item(19, 20).Value = 1.67
$ws.Cells.Item(19, 21).Value = 2.36
$ws.Cells.Item(19, 24).Value = 17.5
$ws.Cells.Item(19, 25).Value = 13
$ws.Cells.Item(19, 26).Value = 19.5
$ws.Cells.Item(19, 28).Value = 13.5
$ws.Cells.Item(19, 29).Value = 8.6
$ws.Cells.Item(19, 30).Value = 12.5
$ws.Cells.Item(19, 31).Value = 28
$ws.Cells.Item(19, 33).Value = 13
$ws.Cells.Item(19, 34).Value = 17
$ws.Cells.Item(19, 35).Value = 38
$ws.Cells.Item(19, 36).Value = 48
$ws.Cells.Item(19, 37).Value = 30
$ws.Cells.Item(19, 39).Value = 1000
$ws.Cells.Item(19, 40).Value = 22
$ws.Cells.Item(19, 41).Value = 21
